$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 20150809
$ws.Range("B4").Value = "perssocialpsychrev"
$ws.Range("C4").Value = "http://psr.sagepub.com/content/by/year/"
$ws.Range("D4").Value = "http://psr.sagepub.com/content/by/year/[0-9]{4}"
$ws.Range("E4").Value = "http://psr.sagepub.com/content/vol[0-9]{1,}/issue[0-9]{1,}/"
$ws.Range("F4").Value = "http://psr.sagepub.com/content/[0-9]{1,}/[0-9]{1,}/[0-9]{1,}.abstract"

$ws.Range("A5").Select()
